$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estonia Meistriliiga")

# Row 32
$ws.Cells.Item(32, 2).Value = 6478314
$ws.Cells.Item(32, 6).Value = 'JK Tammeka Tartu'
$ws.Cells.Item(32, 7).Value = 'JK Tallinna Kalev'
$ws.Cells.Item(32, 8).Value = 1
$ws.Cells.Item(32, 9).Value = 2
$ws.Cells.Item(32, 10).Value = 'A'
$ws.Cells.Item(32, 11).Value = 2.4
$ws.Cells.Item(32, 12).Value = 3.6
$ws.Cells.Item(32, 13).Value = 2.4
$ws.Cells.Item(32, 14).Value = 2.3
$ws.Cells.Item(32, 15).Value = 3.4
$ws.Cells.Item(32, 16).Value = 2.6
$ws.Cells.Item(32, 17).Value = 0
$ws.Cells.Item(32, 18).Value = 1.725
$ws.Cells.Item(32, 19).Value = 1.975
$ws.Cells.Item(32, 20).Value = 2.25
$ws.Cells.Item(32, 21).Value = 1.9
$ws.Cells.Item(32, 22).Value = 1.9
$ws.Cells.Item(32, 23).Value = -1
$ws.Cells.Item(32, 24).Value = -1
$ws.Cells.Item(32, 25).Value = 1.6
$ws.Cells.Item(32, 26).Value = -1
$ws.Cells.Item(32, 27).Value = 0.9750000000000001
$ws.Cells.Item(32, 28).Value = 0.8999999999999999
$ws.Cells.Item(32, 29).Value = -1

# Row 33
$ws.Cells.Item(33, 2).Value = 6475429
$ws.Cells.Item(33, 6).Value = 'Parnu JK Vaprus'
$ws.Cells.Item(33, 7).Value = 'JK Trans Narva'
$ws.Cells.Item(33, 8).Value = 1
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 'H'
$ws.Cells.Item(33, 11).Value = 2.4
$ws.Cells.Item(33, 12).Value = 3.6
$ws.Cells.Item(33, 13).Value = 2.4
$ws.Cells.Item(33, 14).Value = 2.8
$ws.Cells.Item(33, 15).Value = 3.6
$ws.Cells.Item(33, 16).Value = 2.1
$ws.Cells.Item(33, 17).Value = 0.25
$ws.Cells.Item(33, 18).Value = 1.9
$ws.Cells.Item(33, 19).Value = 1.9
$ws.Cells.Item(33, 20).Value = 2.5
$ws.Cells.Item(33, 21).Value = 1.925
$ws.Cells.Item(33, 22).Value = 1.875
$ws.Cells.Item(33, 23).Value = 1.8
$ws.Cells.Item(33, 24).Value = -1
$ws.Cells.Item(33, 25).Value = -1
$ws.Cells.Item(33, 26).Value = 0.8999999999999999
$ws.Cells.Item(33, 27).Value = -1
$ws.Cells.Item(33, 28).Value = -1
$ws.Cells.Item(33, 29).Value = 0.875

# Row 77
$ws.Cells.Item(77, 2).Value = 6139017
$ws.Cells.Item(77, 6).Value = 'JK Tammeka Tartu'
$ws.Cells.Item(77, 7).Value = 'Harju JK Laagri'
$ws.Cells.Item(77, 8).Value = 2
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 'H'
$ws.Cells.Item(77, 11).Value = 1.666
$ws.Cells.Item(77, 12).Value = 3.6
$ws.Cells.Item(77, 13).Value = 4.2
$ws.Cells.Item(77, 14).Value = 1.727
$ws.Cells.Item(77, 15).Value = 3.5
$ws.Cells.Item(77, 16).Value = 4
$ws.Cells.Item(77, 17).Value = -0.75
$ws.Cells.Item(77, 18).Value = 2
$ws.Cells.Item(77, 19).Value = 1.8
$ws.Cells.Item(77, 20).Value = 2.5
$ws.Cells.Item(77, 21).Value = 1.9
$ws.Cells.Item(77, 22).Value = 1.9
$ws.Cells.Item(77, 23).Value = 0.7270000000000001
$ws.Cells.Item(77, 24).Value = -1
$ws.Cells.Item(77, 25).Value = -1
$ws.Cells.Item(77, 26).Value = 1
$ws.Cells.Item(77, 27).Value = -1
$ws.Cells.Item(77, 28).Value = -1
$ws.Cells.Item(77, 29).Value = 0.8999999999999999

# Row 78
$ws.Cells.Item(78, 2).Value = 6139018
$ws.Cells.Item(78, 6).Value = 'JK Tallinna Kalev'
$ws.Cells.Item(78, 7).Value = 'JK Trans Narva'
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 9).Value = 1
$ws.Cells.Item(78, 10).Value = 'A'
$ws.Cells.Item(78, 11).Value = 2.4
$ws.Cells.Item(78, 12).Value = 3.4
$ws.Cells.Item(78, 13).Value = 2.5
$ws.Cells.Item(78, 14).Value = 2.875
$ws.Cells.Item(78, 15).Value = 3.1
$ws.Cells.Item(78, 16).Value = 2.3
$ws.Cells.Item(78, 17).Value = 0.25
$ws.Cells.Item(78, 18).Value = 1.75
$ws.Cells.Item(78, 19).Value = 2.05
$ws.Cells.Item(78, 20).Value = 2.25
$ws.Cells.Item(78, 21).Value = 1.925
$ws.Cells.Item(78, 22).Value = 1.875
$ws.Cells.Item(78, 23).Value = -1
$ws.Cells.Item(78, 24).Value = -1
$ws.Cells.Item(78, 25).Value = 1.3
$ws.Cells.Item(78, 26).Value = -1
$ws.Cells.Item(78, 27).Value = 1.05
$ws.Cells.Item(78, 28).Value = -1
$ws.Cells.Item(78, 29).Value = 0.875

# Row 177
$ws.Cells.Item(177, 2).Value = 6537957
$ws.Cells.Item(177, 6).Value = 'FC Flora Tallinn'
$ws.Cells.Item(177, 7).Value = 'JK Nomme Kalju'
$ws.Cells.Item(177, 8).Value = 0
$ws.Cells.Item(177, 9).Value = 0
$ws.Cells.Item(177, 10).Value = 'D'
$ws.Cells.Item(177, 11).Value = 1.4
$ws.Cells.Item(177, 12).Value = 4
$ws.Cells.Item(177, 13).Value = 7.5
$ws.Cells.Item(177, 14).Value = 1.5
$ws.Cells.Item(177, 15).Value = 4.2
$ws.Cells.Item(177, 16).Value = 5
$ws.Cells.Item(177, 17).Value = -1
$ws.Cells.Item(177, 18).Value = 1.85
$ws.Cells.Item(177, 19).Value = 1.95
$ws.Cells.Item(177, 20).Value = 2.75
$ws.Cells.Item(177, 21).Value = 1.85
$ws.Cells.Item(177, 22).Value = 1.95
$ws.Cells.Item(177, 23).Value = -1
$ws.Cells.Item(177, 24).Value = 3.2
$ws.Cells.Item(177, 25).Value = -1
$ws.Cells.Item(177, 26).Value = -1
$ws.Cells.Item(177, 27).Value = 0.95
$ws.Cells.Item(177, 28).Value = -1
$ws.Cells.Item(177, 29).Value = 0.95

# Row 178
$ws.Cells.Item(178, 2).Value = 6537869
$ws.Cells.Item(178, 6).Value = 'JK Tallinna Kalev'
$ws.Cells.Item(178, 7).Value = 'JK Trans Narva'
$ws.Cells.Item(178, 8).Value = 5
$ws.Cells.Item(178, 9).Value = 0
$ws.Cells.Item(178, 10).Value = 'H'
$ws.Cells.Item(178, 11).Value = 1.6
$ws.Cells.Item(178, 12).Value = 4
$ws.Cells.Item(178, 13).Value = 4.5
$ws.Cells.Item(178, 14).Value = 1.65
$ws.Cells.Item(178, 15).Value = 4
$ws.Cells.Item(178, 16).Value = 4.333
$ws.Cells.Item(178, 17).Value = -0.75
$ws.Cells.Item(178, 18).Value = 1.8
$ws.Cells.Item(178, 19).Value = 2
$ws.Cells.Item(178, 20).Value = 2.75
$ws.Cells.Item(178, 21).Value = 1.9
$ws.Cells.Item(178, 22).Value = 1.9
$ws.Cells.Item(178, 23).Value = 0.6499999999999999
$ws.Cells.Item(178, 24).Value = -1
$ws.Cells.Item(178, 25).Value = -1
$ws.Cells.Item(178, 26).Value = 0.8
$ws.Cells.Item(178, 27).Value = -1
$ws.Cells.Item(178, 28).Value = 0.8999999999999999
$ws.Cells.Item(178, 29).Value = -1

# Row 179
$ws.Cells.Item(179, 2).Value = 6535416
$ws.Cells.Item(179, 6).Value = 'Paide Linnameeskond'
$ws.Cells.Item(179, 7).Value = 'FC Levadia Tallinn'
$ws.Cells.Item(179, 8).Value = 2
$ws.Cells.Item(179, 9).Value = 2
$ws.Cells.Item(179, 10).Value = 'D'
$ws.Cells.Item(179, 11).Value = 3
$ws.Cells.Item(179, 12).Value = 3.8
$ws.Cells.Item(179, 13).Value = 2
$ws.Cells.Item(179, 14).Value = 3
$ws.Cells.Item(179, 15).Value = 4
$ws.Cells.Item(179, 16).Value = 1.909
$ws.Cells.Item(179, 17).Value = 0.5
$ws.Cells.Item(179, 18).Value = 1.85
$ws.Cells.Item(179, 19).Value = 1.95
$ws.Cells.Item(179, 20).Value = 2.75
$ws.Cells.Item(179, 21).Value = 1.95
$ws.Cells.Item(179, 22).Value = 1.85
$ws.Cells.Item(179, 23).Value = -1
$ws.Cells.Item(179, 24).Value = 3
$ws.Cells.Item(179, 25).Value = -1
$ws.Cells.Item(179, 26).Value = 0.8500000000000001
$ws.Cells.Item(179, 27).Value = -1
$ws.Cells.Item(179, 28).Value = 0.95
$ws.Cells.Item(179, 29).Value = -1

# Row 180
$ws.Cells.Item(180, 2).Value = 6533597
$ws.Cells.Item(180, 6).Value = 'FC Kuressaare'
$ws.Cells.Item(180, 7).Value = 'Parnu JK Vaprus'
$ws.Cells.Item(180, 8).Value = 1
$ws.Cells.Item(180, 9).Value = 0
$ws.Cells.Item(180, 10).Value = 'H'
$ws.Cells.Item(180, 11).Value = 2.5
$ws.Cells.Item(180, 12).Value = 3.4
$ws.Cells.Item(180, 13).Value = 2.5
$ws.Cells.Item(180, 14).Value = 2.15
$ws.Cells.Item(180, 15).Value = 3.6
$ws.Cells.Item(180, 16).Value = 2.875
$ws.Cells.Item(180, 17).Value = -0.25
$ws.Cells.Item(180, 18).Value = 1.95
$ws.Cells.Item(180, 19).Value = 1.85
$ws.Cells.Item(180, 20).Value = 2.75
$ws.Cells.Item(180, 21).Value = 1.95
$ws.Cells.Item(180, 22).Value = 1.85
$ws.Cells.Item(180, 23).Value = 1.15
$ws.Cells.Item(180, 24).Value = -1
$ws.Cells.Item(180, 25).Value = -1
$ws.Cells.Item(180, 26).Value = 0.95
$ws.Cells.Item(180, 27).Value = -1
$ws.Cells.Item(180, 28).Value = -1
$ws.Cells.Item(180, 29).Value = 0.8500000000000001
